$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the numeric values in B2:D4 to the new data
$ws.Range("B2").Value = 0.93735980305625011
$ws.Range("C2").Value = 0.27896708354690442
$ws.Range("D2").Value = -0.20864794729803221

$ws.Range("B3").Value = 0.33714882217054887
$ws.Range("C3").Value = -0.57572686584418165
$ws.Range("D3").Value = 0.74489546088981329

$ws.Range("B4").Value = 0.087677085509040339
$ws.Range("C4").Value = -0.76858047219700454
$ws.Range("D4").Value = -0.63371664522408611

# Match the saved selection/view state (whole data range selected)
$ws.Range("A1:D4").Select()

# Remove the extra (now unused) worksheets, keeping only Sheet1
$wb.Worksheets.Item("Sheet3").Delete()
$wb.Worksheets.Item("Sheet2").Delete()
